$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rose")

# Update a couple of "Gol" (goals) values for existing players
$ws.Range("E12").Value = 1
$ws.Range("E15").Value = 1

# Swap jersey numbers between two existing "PM Sport" players
$ws.Range("D65").Value = 6
$ws.Range("D66").Value = 5

# Insert a new roster row for a newly added player "Maicol Batti" (PM Sport,
# Attaccante, #7) above the current row 72, shifting the remaining rows down
$ws.Rows.Item(72).Insert()
$ws.Range("A72").Value = "PM Sport"
$ws.Range("B72").Value = "Maicol Batti"
$ws.Range("C72").Value = "Attaccante"
$ws.Range("D72").Value = 7
$ws.Range("E72").Value = 0

# Widen column F slightly (matches the author's workbook-view tweak)
$ws.Columns.Item(6).ColumnWidth = 16.83

# Restore the active selection used when the file was saved
$ws.Range("B14").Select()
